$d = $word.ActiveDocument

$replacements = @(
    @("69÷8=", "85÷5="),
    @("65÷2=", "97÷3="),
    @("43÷6=", "10÷9="),
    @("29÷6=", "92÷8="),
    @("98÷6=", "37÷8="),
    @("37÷6=", "61÷8="),
    @("22÷6=", "79÷2="),
    @("13÷6=", "32÷2="),
    @("18÷8=", "82÷2="),
    @("52÷4=", "61÷9="),
    @("42÷9=", "24÷3="),
    @("60÷5=", "11÷2="),
    @("20÷3=", "94÷2="),
    @("47÷9=", "33÷5="),
    @("47÷6=", "64÷7="),
    @("58÷6=", "72÷3="),
    @("24÷9=", "63÷9="),
    @("69÷3=", "99÷7="),
    @("34÷2=", "87÷6="),
    @("62÷4=", "46÷6="),
    @("25÷9=", "38÷4="),
    @("74÷7=", "68÷4="),
    @("44÷6=", "25÷6="),
    @("40÷6=", "33÷2="),
    @("21÷7=", "91÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
